$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.54861307144165
$ws.Range("B1").Value = 3.045896530151367
$ws.Range("C1").Value = 0.6381919384002686
$ws.Range("D1").Value = 0.6112232804298401
$ws.Range("E1").Value = 0.2427859902381897
